$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 470
$ws.Range("I6").Value = 470
$ws.Range("K6").Value = 1410
$ws.Range("M6").Value = -1298

$ws.Range("H8").Value = 8759.125
$ws.Range("I8").Value = 8759.125
$ws.Range("K8").Value = 26277.375
$ws.Range("M8").Value = -26138.375

$ws.Range("H17").Value = 908.9231
$ws.Range("J17").Value = 908.9231
$ws.Range("L17").Value = 2726.7693
$ws.Range("N17").Value = -3062.7693

$ws.Range("H28").Value = 1719.4117
$ws.Range("I28").Value = 1172.6923
$ws.Range("J28").Value = 3496.25
$ws.Range("K28").Value = 1172.6923
$ws.Range("L28").Value = 3496.25
$ws.Range("M28").Value = -687.6922999999999
$ws.Range("N28").Value = -4466.25

$ws.Range("H135").Value = 28572624
$ws.Range("I135").Value = 1120.9615
$ws.Range("J135").Value = 111112520
$ws.Range("K135").Value = 10088.6535
$ws.Range("L135").Value = 1000012680
$ws.Range("M135").Value = -7553.653499999999
$ws.Range("N135").Value = -1000017750

$ws.Range("H137").Value = 1670.6364
$ws.Range("I137").Value = 1274.1
$ws.Range("J137").Value = 2280.6924
$ws.Range("K137").Value = 3822.3
$ws.Range("L137").Value = 6842.0772
$ws.Range("M137").Value = -1272.3
$ws.Range("N137").Value = -11942.0772

$ws.Range("H138").Value = 4764833
$ws.Range("I138").Value = 2334
$ws.Range("J138").Value = 6669832.5
$ws.Range("K138").Value = 7002
$ws.Range("L138").Value = 20009497.5
$ws.Range("M138").Value = -1862
$ws.Range("N138").Value = -20019777.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1350
$ws.Range("I41").Value = 1350
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1350
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -936
$ws.Range("N41").ClearContents()

$ws.Range("H63").Value = 3887.75
$ws.Range("I63").Value = 2925.5
$ws.Range("J63").Value = 4850
$ws.Range("K63").Value = 2925.5
$ws.Range("L63").Value = 4850
$ws.Range("M63").Value = -2239.5
$ws.Range("N63").Value = -6222

$ws.Range("H66").Value = 3887.75
$ws.Range("I66").Value = 2925.5
$ws.Range("J66").Value = 4850
$ws.Range("K66").Value = 14627.5
$ws.Range("L66").Value = 24250
$ws.Range("M66").Value = -11195.5
$ws.Range("N66").Value = -31114

$ws.Range("H74").Value = 29196.25
$ws.Range("I74").Value = 48899.285
$ws.Range("J74").Value = 1612
$ws.Range("K74").Value = 48899.285
$ws.Range("L74").Value = 1612
$ws.Range("M74").Value = -48025.285
$ws.Range("N74").Value = -3360

$ws.Range("H77").Value = 29196.25
$ws.Range("I77").Value = 48899.285
$ws.Range("J77").Value = 1612
$ws.Range("K77").Value = 244496.425
$ws.Range("L77").Value = 8060
$ws.Range("M77").Value = -240128.425
$ws.Range("N77").Value = -16796

$ws.Range("H138").Value = 36999
$ws.Range("J138").Value = 36999
$ws.Range("L138").Value = 36999
$ws.Range("N138").Value = -47279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 528987.0600000001
$ws.Range("I134").Value = 802522.9
$ws.Range("J134").Value = 2956.6538
$ws.Range("K134").Value = 2407568.7
$ws.Range("L134").Value = 8869.9614
$ws.Range("M134").Value = -2405033.7
$ws.Range("N134").Value = -13939.9614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1467.5438
$ws.Range("I31").Value = 897.75
$ws.Range("J31").Value = 2444.3333
$ws.Range("K31").Value = 897.75
$ws.Range("L31").Value = 2444.3333
$ws.Range("M31").Value = -602.75
$ws.Range("N31").Value = -3034.3333

$ws.Range("H34").Value = 1467.5438
$ws.Range("I34").Value = 897.75
$ws.Range("J34").Value = 2444.3333
$ws.Range("K34").Value = 897.75
$ws.Range("L34").Value = 2444.3333
$ws.Range("M34").Value = -695.75
$ws.Range("N34").Value = -2848.3333

$ws.Range("H58").Value = 3144.9788
$ws.Range("I58").Value = 3355.6191
$ws.Range("J58").Value = 1375.6
$ws.Range("K58").Value = 3355.6191
$ws.Range("L58").Value = 1375.6
$ws.Range("M58").Value = -3152.6191
$ws.Range("N58").Value = -1781.6

$ws.Range("H132").Value = 951695.4
$ws.Range("I132").Value = 1835.9395
$ws.Range("J132").Value = 6175922.5
$ws.Range("K132").Value = 5507.818499999999
$ws.Range("L132").Value = 18527767.5
$ws.Range("M132").Value = -2977.818499999999
$ws.Range("N132").Value = -18532827.5

$ws.Range("H134").Value = 1982
$ws.Range("I134").Value = 1924.2572
$ws.Range("J134").Value = 2318.8333
$ws.Range("K134").Value = 5772.7716
$ws.Range("L134").Value = 6956.499899999999
$ws.Range("M134").Value = -3237.7716
$ws.Range("N134").Value = -12026.4999

$ws.Range("H136").Value = 3144.9788
$ws.Range("I136").Value = 3355.6191
$ws.Range("J136").Value = 1375.6
$ws.Range("K136").Value = 10066.8573
$ws.Range("L136").Value = 4126.799999999999
$ws.Range("M136").Value = -7516.8573
$ws.Range("N136").Value = -9226.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 638.3333
$ws.Range("J80").Value = 638.3333
$ws.Range("L80").Value = 1914.9999
$ws.Range("N80").Value = -3786.9999

$ws.Range("H83").Value = 638.3333
$ws.Range("J83").Value = 638.3333
$ws.Range("L83").Value = 5744.9997
$ws.Range("N83").Value = -15104.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 37375.2
$ws.Range("I122").Value = 46635.695
$ws.Range("J122").Value = 6947.857
$ws.Range("K122").Value = 139907.085
$ws.Range("L122").Value = 20843.571
$ws.Range("M122").Value = -137457.085
$ws.Range("N122").Value = -25743.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 857.8929000000001
$ws.Range("I46").Value = 1125.8572
$ws.Range("J46").Value = 768.5714
$ws.Range("K46").Value = 1125.8572
$ws.Range("L46").Value = 768.5714
$ws.Range("M46").Value = -937.8571999999999
$ws.Range("N46").Value = -1144.5714

$ws.Range("H93").Value = 39732.223
$ws.Range("I93").Value = 897.7
$ws.Range("J93").Value = 88275.375
$ws.Range("K93").Value = 897.7
$ws.Range("L93").Value = 88275.375
$ws.Range("M93").Value = 350.3
$ws.Range("N93").Value = -90771.375

$ws.Range("H136").Value = 1398.2037
$ws.Range("I136").Value = 919.8
$ws.Range("J136").Value = 2279.4736
$ws.Range("K136").Value = 2759.4
$ws.Range("L136").Value = 6838.4208
$ws.Range("M136").Value = -209.3999999999996
$ws.Range("N136").Value = -11938.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 29111.889
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 29111.889
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 29111.889
$ws.Range("N12").Value = -29395.889
$ws.Range("M12").ClearContents()

$ws.Range("H136").Value = 1050.0308
$ws.Range("I136").Value = 650.02563
$ws.Range("J136").Value = 1650.0385
$ws.Range("K136").Value = 1950.07689
$ws.Range("L136").Value = 4950.1155
$ws.Range("M136").Value = 599.9231100000002
$ws.Range("N136").Value = -10050.1155
